$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header values: A1=Ano, B1=Cadastrado, C1=Sem Cadastro
$ws.Range("A1").Value = "Ano"
$ws.Range("B1").Value = "Cadastrado"
$ws.Range("C1").Value = "Sem Cadastro"

# New data set (years 2018-2025)
$data = @(
    @(2018, 14294, 685903.33),
    @(2019, 128163.6, 1613553.48),
    @(2020, 543045.64, 2770811.39),
    @(2021, 1456784.01, 5107712.11),
    @(2022, 1655989.08, 5704577.96),
    @(2023, 1148594.77, 5784239.8),
    @(2024, 2797633.46, 6627864.05),
    @(2025, 1010657.76, 2254826.4)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Remove the bold/bordered style previously applied to column A data cells (A2:A6)
# so the new data rows (A2:A9) have the default style
$ws.Range("A2:A9").Style = "Normal"
